$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data (SKU/Qty/Cost columns) is stored as text in this sheet,
# even though many values look numeric. Force the new cells to Text format
# first so typed-looking values ("150300865", "79.99", ...) are written as
# strings instead of being auto-converted to numbers.
$newRange = $ws.Range("A23:E24")
$newRange.NumberFormat = "@"

$ws.Range("A23").Value = "150300865"
$ws.Range("B23").Value = "Bag Paper - 6x13.5 Window"
$ws.Range("C23").Value = "4"
$ws.Range("D23").Value = "79.99"
$ws.Range("E23").Value = "319.96"

$ws.Range("A24").Value = "588MILK632"
$ws.Range("B24").Value = "Urnex - Rinza"
$ws.Range("C24").Value = "24"
$ws.Range("D24").Value = "17.99"
$ws.Range("E24").Value = "431.76"
